# Auto-generated edit script applying numeric corrections
# to the Jenova_Profits workbook per the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 49712.535
$ws.Range("J17").Value = 49712.535
$ws.Range("L17").Value = 149137.605
$ws.Range("N17").Value = -149473.605
$ws.Range("H32").Value = 2418.3333
$ws.Range("I32").Value = 2422.8572
$ws.Range("J32").Value = 2416.0715
$ws.Range("K32").Value = 2422.8572
$ws.Range("L32").Value = 2416.0715
$ws.Range("M32").Value = -2096.8572
$ws.Range("N32").Value = -3068.0715
$ws.Range("H64").Value = 6271.5
$ws.Range("I64").Value = 5073.9165
$ws.Range("J64").Value = 8666.666999999999
$ws.Range("K64").Value = 5073.9165
$ws.Range("L64").Value = 8666.666999999999
$ws.Range("M64").Value = -4825.9165
$ws.Range("N64").Value = -9162.666999999999
$ws.Range("H67").Value = 6271.5
$ws.Range("I67").Value = 5073.9165
$ws.Range("J67").Value = 8666.666999999999
$ws.Range("K67").Value = 5073.9165
$ws.Range("L67").Value = 8666.666999999999
$ws.Range("M67").Value = -4215.9165
$ws.Range("N67").Value = -10382.667
$ws.Range("H70").Value = 114194.78
$ws.Range("I70").Value = 1750
$ws.Range("J70").Value = 146321.86
$ws.Range("K70").Value = 5250
$ws.Range("L70").Value = 438965.58
$ws.Range("M70").Value = -4980
$ws.Range("N70").Value = -439505.58
$ws.Range("H73").Value = 114194.78
$ws.Range("I73").Value = 1750
$ws.Range("J73").Value = 146321.86
$ws.Range("K73").Value = 5250
$ws.Range("L73").Value = 438965.58
$ws.Range("M73").Value = -4314
$ws.Range("N73").Value = -440837.58
$ws.Range("H113").Value = 3373.6
$ws.Range("J113").Value = 2992
$ws.Range("L113").Value = 2992
$ws.Range("N113").Value = -9500
$ws.Range("H131").Value = 2269.8965
$ws.Range("I131").Value = 1519.5217
$ws.Range("K131").Value = 4558.5651
$ws.Range("M131").Value = 481.4349000000002
$ws.Range("H135").Value = 2378.625
$ws.Range("I135").Value = 932.7143
$ws.Range("K135").Value = 8394.4287
$ws.Range("M135").Value = -5859.4287
$ws.Range("H138").Value = 5344.347
$ws.Range("I138").Value = 2879.375
$ws.Range("J138").Value = 5825.317
$ws.Range("K138").Value = 8638.125
$ws.Range("L138").Value = 17475.951
$ws.Range("M138").Value = -3498.125
$ws.Range("N138").Value = -27755.951
$ws.Range("H139").Value = 50000
$ws.Range("J139").Value = 50000
$ws.Range("L139").Value = 50000
$ws.Range("N139").Value = -60280
$ws.Range("H141").Value = 4934.769
$ws.Range("J141").Value = 4053.5
$ws.Range("L141").Value = 12160.5
$ws.Range("N141").Value = -22520.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2314.8333
$ws.Range("I32").Value = 2015.4783
$ws.Range("K32").Value = 2015.4783
$ws.Range("M32").Value = -1728.4783
$ws.Range("H61").Value = 2175.1035
$ws.Range("I61").Value = 1489.5652
$ws.Range("K61").Value = 1489.5652
$ws.Range("M61").Value = -1277.5652
$ws.Range("H74").Value = 17103.277
$ws.Range("I74").Value = 25582.545
$ws.Range("K74").Value = 25582.545
$ws.Range("M74").Value = -24708.545
$ws.Range("H77").Value = 17103.277
$ws.Range("I77").Value = 25582.545
$ws.Range("K77").Value = 127912.725
$ws.Range("M77").Value = -123544.725
$ws.Range("H136").Value = 2175.1035
$ws.Range("I136").Value = 1489.5652
$ws.Range("K136").Value = 4468.6956
$ws.Range("M136").Value = -1918.6956
$ws.Range("H140").Value = 49800
$ws.Range("J140").Value = 49800
$ws.Range("L140").Value = 49800
$ws.Range("N140").Value = -60160

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H87").Value = 45000
$ws.Range("I87").Value = 45000
$ws.Range("K87").Value = 45000
$ws.Range("M87").Value = -43752
$ws.Range("H90").Value = 45000
$ws.Range("I90").Value = 45000
$ws.Range("K90").Value = 135000
$ws.Range("M90").Value = -128760
$ws.Range("H102").Value = 5541.6665
$ws.Range("I102").Value = 5541.6665
$ws.Range("K102").Value = 5541.6665
$ws.Range("M102").Value = -2296.6665
$ws.Range("H134").Value = 3927.64
$ws.Range("I134").Value = 2504.4443
$ws.Range("K134").Value = 7513.3329
$ws.Range("M134").Value = -4978.3329

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3030.5625
$ws.Range("I31").Value = 1653.8182
$ws.Range("K31").Value = 1653.8182
$ws.Range("M31").Value = -1358.8182
$ws.Range("H34").Value = 3030.5625
$ws.Range("I34").Value = 1653.8182
$ws.Range("K34").Value = 1653.8182
$ws.Range("M34").Value = -1451.8182
$ws.Range("H58").Value = 325437.8
$ws.Range("I58").Value = 527575.9
$ws.Range("J58").Value = 5385.8335
$ws.Range("K58").Value = 527575.9
$ws.Range("L58").Value = 5385.8335
$ws.Range("M58").Value = -527372.9
$ws.Range("N58").Value = -5791.8335
$ws.Range("H60").Value = 73333.336
$ws.Range("J60").Value = 73333.336
$ws.Range("L60").Value = 73333.336
$ws.Range("N60").Value = -74355.336
$ws.Range("H62").Value = 3325
$ws.Range("I62").Value = 3348.125
$ws.Range("J62").Value = 3294.1667
$ws.Range("K62").Value = 3348.125
$ws.Range("L62").Value = 3294.1667
$ws.Range("M62").Value = -2724.125
$ws.Range("N62").Value = -4542.1667
$ws.Range("H65").Value = 3325
$ws.Range("I65").Value = 3348.125
$ws.Range("J65").Value = 3294.1667
$ws.Range("K65").Value = 16740.625
$ws.Range("L65").Value = 16470.8335
$ws.Range("M65").Value = -13620.625
$ws.Range("N65").Value = -22710.8335
$ws.Range("H132").Value = 2643.3333
$ws.Range("I132").Value = 2366.9167
$ws.Range("K132").Value = 7100.750100000001
$ws.Range("M132").Value = -4570.750100000001
$ws.Range("H136").Value = 325437.8
$ws.Range("I136").Value = 527575.9
$ws.Range("J136").Value = 5385.8335
$ws.Range("K136").Value = 1582727.7
$ws.Range("L136").Value = 16157.5005
$ws.Range("M136").Value = -1580177.7
$ws.Range("N136").Value = -21257.5005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 1987.5
$ws.Range("I46").Value = 1975
$ws.Range("K46").Value = 5925
$ws.Range("M46").Value = -5834

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 25371.875
$ws.Range("I57").Value = 13829.167
$ws.Range("K57").Value = 13829.167
$ws.Range("M57").Value = -13009.167
$ws.Range("H80").Value = 1669762.6
$ws.Range("I80").Value = 1669255.6
$ws.Range("K80").Value = 1669255.6
$ws.Range("M80").Value = -1668257.6
$ws.Range("H83").Value = 1669762.6
$ws.Range("I83").Value = 1669255.6
$ws.Range("K83").Value = 8346278
$ws.Range("M83").Value = -8341286
$ws.Range("H132").Value = 273515.7
$ws.Range("I132").Value = 315408.8
$ws.Range("K132").Value = 946226.3999999999
$ws.Range("M132").Value = -943696.3999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 679.3
$ws.Range("I22").Value = 667.25
$ws.Range("K22").Value = 667.25
$ws.Range("M22").Value = -372.25
$ws.Range("H27").Value = 679.3
$ws.Range("I27").Value = 667.25
$ws.Range("K27").Value = 667.25
$ws.Range("M27").Value = -560.25
$ws.Range("H40").Value = 912697.0600000001
$ws.Range("I40").Value = 912697.0600000001
$ws.Range("K40").Value = 912697.0600000001
$ws.Range("M40").Value = -912561.0600000001
$ws.Range("H55").Value = 790.1
$ws.Range("I55").Value = 600.125
$ws.Range("K55").Value = 600.125
$ws.Range("M55").Value = -427.125
$ws.Range("H132").Value = 6443.8887
$ws.Range("I132").Value = 4998
$ws.Range("K132").Value = 14994
$ws.Range("M132").Value = -12464
$ws.Range("H136").Value = 3516.4443
$ws.Range("I136").Value = 3027.2
$ws.Range("K136").Value = 9081.599999999999
$ws.Range("M136").Value = -6531.599999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2876.1428
$ws.Range("J81").Value = 4407.3335
$ws.Range("L81").Value = 8814.666999999999
$ws.Range("N81").Value = -10936.667
$ws.Range("H84").Value = 2876.1428
$ws.Range("J84").Value = 4407.3335
$ws.Range("L84").Value = 44073.335
$ws.Range("N84").Value = -54681.335
$ws.Range("H101").Value = 30000
$ws.Range("J101").Value = 30000
$ws.Range("L101").Value = 30000
$ws.Range("N101").Value = -36490
$ws.Range("H133").Value = 55798.6
$ws.Range("J133").Value = 55798.6
$ws.Range("L133").Value = 55798.6
$ws.Range("N133").Value = -65918.60000000001
$ws.Range("H136").Value = 528488.9399999999
$ws.Range("I136").Value = 528488.9399999999
$ws.Range("K136").Value = 1585466.82
$ws.Range("M136").Value = -1582916.82

Write-Host "Applied 214 cell updates"